$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("gamma")
$ws.Name = "summary"
